$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: row 2 already exists (will become row with id 165/fra),
# insert 9 more rows below it so we end up with rows 2..11 for data ---
$ws.Range("A3:A11").EntireRow.Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = "lang_code"
$ws.Cells.Item(1,2).Value = "id"
$ws.Cells.Item(1,3).Value = "name"
$ws.Cells.Item(1,4).Value = "brand"
$ws.Cells.Item(1,5).Value = "model"
$ws.Cells.Item(1,6).Value = "dtyp_code"
$ws.Cells.Item(1,7).Value = "min_driver_ver"
$ws.Cells.Item(1,8).Value = "descr"
$ws.Cells.Item(1,9).Value = "is_active"

# --- Data rows 2..11 ---
# row, lang_code, id, name, brand, model, dtyp_code, min_driver_ver, descr, is_active
$rows = @(
  @(2,  "fra", 165, "Fingerprint Scanner",          "Safran Morpho", "1300 E2",  "FRS", 1.12,  "To scan fingerprint", $true),
  @(3,  "eng", 165, "Fingerprint Scanner",          "Safran Morpho", "1300 E2",  "FRS", 1.12,  "To scan fingerprint", $true),
  @(4,  "fra", 327, "High Speed Dual Iris Scanner", "Cogent",        "3M",       "IRS", 2.34,  "To scan iris",        $true),
  @(5,  "eng", 327, "High Speed Dual Iris Scanner", "Cogent",        "3M",       "IRS", 2.34,  "To scan iris",        $true),
  @(6,  "fra", 736, "Webcam",                       "Logitech",      "C270",     "CMR", 2.086, "To capture photo",    $true),
  @(7,  "eng", 736, "Webcam",                       "Logitech",      "C270",     "CMR", 2.086, "To capture photo",    $true),
  @(8,  "fra", 801, "imageFORMULA",                 "Canon",         "DR-C130",  "SCN", 1.02,  "To scan documents",   $true),
  @(9,  "eng", 801, "imageFORMULA",                 "Canon",         "DR-C130",  "SCN", 1.02,  "To scan documents",   $true),
  @(10, "fra", 920, "Single Function Inkjet",       "Canon",         "TS207",    "PRT", 1.123, "To print documents",  $true),
  @(11, "eng", 920, "Single Function Inkjet",       "Canon",         "TS207",    "PRT", 1.123, "To print documents",  $true)
)

foreach ($r in $rows) {
    $rowIdx = $r[0]
    $ws.Cells.Item($rowIdx, 1).Value = $r[1]
    $ws.Cells.Item($rowIdx, 2).Value = $r[2]
    $ws.Cells.Item($rowIdx, 3).Value = $r[3]
    $ws.Cells.Item($rowIdx, 4).Value = $r[4]
    $ws.Cells.Item($rowIdx, 5).Value = $r[5]
    $ws.Cells.Item($rowIdx, 6).Value = $r[6]
    $ws.Cells.Item($rowIdx, 7).Value = $r[7]
    $ws.Cells.Item($rowIdx, 8).Value = $r[8]
    $ws.Cells.Item($rowIdx, 9).Value = $r[9]
}

# --- Column A (lang_code) data cells use the same bold/bordered/centered
# style as the header row, matching the source file ---
$ws.Range("A1").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
